$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 159: 14.03.2024 update
$ws.Range("A159").Value = "14.03.2024"
$ws.Range("B159").Value = "14.03.2024"
$ws.Range("C159").Value = 31341
$ws.Range("D159").Value = 12300
$ws.Range("E159").Value = 8400
$ws.Range("F159").Value = 73134
$ws.Range("G159").Value = 8663
$ws.Range("H159").Value = 6327
$ws.Range("I159").Value = 8000
$ws.Range("J159").Value = 433
$ws.Range("K159").Value = 116
$ws.Range("L159").Value = 4650
$ws.Range("M159").Value = "https://web.archive.org/web/20240314205339/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"

# Row 160: 15.03.2024 update
$ws.Range("A160").Value = "15.03.2024"
$ws.Range("B160").Value = "15.03.2024"
$ws.Range("C160").Value = 31490
$ws.Range("D160").Value = 12300
$ws.Range("E160").Value = 8400
$ws.Range("F160").Value = 73439
$ws.Range("G160").Value = 8663
$ws.Range("H160").Value = 6327
$ws.Range("I160").Value = 8000
$ws.Range("J160").Value = 433
$ws.Range("K160").Value = 116
$ws.Range("L160").Value = 4650
$ws.Range("M160").Value = "https://web.archive.org/web/20240315162911/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"
